# Update the LR-pairs_lrc2p sheet (Il1rn-Il1r1) with newly recomputed TPM-based
# ligand/receptor expression, specificity, and edge-weight statistics.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.098365
$ws.Range("H2").Value = 0.295095
$ws.Range("I2").Value = 0.001274112831990382
$ws.Range("J2").Value = 0.001274112831990381
$ws.Range("M2").Value = 12.673913
$ws.Range("N2").Value = 38.021739
$ws.Range("O2").Value = 0.234043494199914
$ws.Range("P2").Value = 0.234043494199914
$ws.Range("Q2").Value = 1.246669452245
$ws.Range("R2").Value = 11.220025070205
$ws.Range("S2").Value = 0.0002981978192039769
$ws.Range("T2").Value = 0.0002981978192039768
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.098365
$ws.Range("H3").Value = 0.295095
$ws.Range("I3").Value = 0.001274112831990382
$ws.Range("J3").Value = 0.001274112831990381
$ws.Range("O3").Value = 0.6388539132363013
$ws.Range("P3").Value = 0.6388539132363011
$ws.Range("Q3").Value = 3.402955765984999
$ws.Range("R3").Value = 30.626601893865
$ws.Range("S3").Value = 0.0008139719686216413
$ws.Range("T3").Value = 0.000813971968621641
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.098365
$ws.Range("H4").Value = 0.295095
$ws.Range("I4").Value = 0.001274112831990382
$ws.Range("J4").Value = 0.001274112831990381
$ws.Range("M4").Value = 6.728406666666667
$ws.Range("N4").Value = 20.18522
$ws.Range("O4").Value = 0.1242504825987572
$ws.Range("P4").Value = 0.1242504825987572
$ws.Range("Q4").Value = 0.6618397217666667
$ws.Range("R4").Value = 5.9565574959
$ws.Range("S4").Value = 0.0001583091342600742
$ws.Range("T4").Value = 0.0001583091342600741
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.098365
$ws.Range("H5").Value = 0.295095
$ws.Range("I5").Value = 0.001274112831990382
$ws.Range("J5").Value = 0.001274112831990381
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.1544473333333333
$ws.Range("N5").Value = 0.463342
$ws.Range("O5").Value = 0.002852109965027549
$ws.Range("P5").Value = 0.002852109965027548
$ws.Range("Q5").Value = 0.01519221194333333
$ws.Range("R5").Value = 0.13672990749
$ws.Range("S5").Value = "0.000003633909904689238"
$ws.Range("T5").Value = "0.000003633909904689237"
$ws.Range("G6").Value = 0.03342533333333333
$ws.Range("H6").Value = 0.100276
$ws.Range("I6").Value = 0.0004329552799629527
$ws.Range("J6").Value = 0.0004329552799629526
$ws.Range("M6").Value = 12.673913
$ws.Range("N6").Value = 38.021739
$ws.Range("O6").Value = 0.234043494199914
$ws.Range("P6").Value = 0.234043494199914
$ws.Range("Q6").Value = 0.4236297666626667
$ws.Range("R6").Value = 3.812667899964
$ws.Range("S6").Value = 0.0001013303665548315
$ws.Range("T6").Value = 0.0001013303665548314
$ws.Range("G7").Value = 0.03342533333333333
$ws.Range("H7").Value = 0.100276
$ws.Range("I7").Value = 0.0004329552799629527
$ws.Range("J7").Value = 0.0004329552799629526
$ws.Range("O7").Value = 0.6388539132363013
$ws.Range("P7").Value = 0.6388539132363011
$ws.Range("Q7").Value = 1.156355724054667
$ws.Range("R7").Value = 10.407201516492
$ws.Range("S7").Value = 0.0002765951748606507
$ws.Range("T7").Value = 0.0002765951748606506
$ws.Range("G8").Value = 0.03342533333333333
$ws.Range("H8").Value = 0.100276
$ws.Range("I8").Value = 0.0004329552799629527
$ws.Range("J8").Value = 0.0004329552799629526
$ws.Range("M8").Value = 6.728406666666667
$ws.Range("N8").Value = 20.18522
$ws.Range("O8").Value = 0.1242504825987572
$ws.Range("P8").Value = 0.1242504825987572
$ws.Range("Q8").Value = 0.2248992356355556
$ws.Range("R8").Value = 2.02409312072
$ws.Range("S8").Value = "0.00005379490247907692"
$ws.Range("T8").Value = "0.0000537949024790769"
$ws.Range("G9").Value = 0.03342533333333333
$ws.Range("H9").Value = 0.100276
$ws.Range("I9").Value = 0.0004329552799629527
$ws.Range("J9").Value = 0.0004329552799629526
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.1544473333333333
$ws.Range("N9").Value = 0.463342
$ws.Range("O9").Value = 0.002852109965027549
$ws.Range("P9").Value = 0.002852109965027548
$ws.Range("Q9").Value = 0.005162453599111111
$ws.Range("R9").Value = 0.046462082392
$ws.Range("S9").Value = "0.00000123483606839363"
$ws.Range("T9").Value = "0.000001234836068393629"
$ws.Range("G10").Value = 77.07094833333333
$ws.Range("H10").Value = 231.212845
$ws.Range("I10").Value = 0.9982929318880467
$ws.Range("J10").Value = 0.9982929318880466
$ws.Range("M10").Value = 12.673913
$ws.Range("N10").Value = 38.021739
$ws.Range("O10").Value = 0.234043494199914
$ws.Range("P10").Value = 0.234043494199914
$ws.Range("Q10").Value = 976.7904940041616
$ws.Range("R10").Value = 8791.114446037454
$ws.Range("S10").Value = 0.2336439660141552
$ws.Range("T10").Value = 0.2336439660141552
$ws.Range("G11").Value = 77.07094833333333
$ws.Range("H11").Value = 231.212845
$ws.Range("I11").Value = 0.9982929318880467
$ws.Range("J11").Value = 0.9982929318880466
$ws.Range("O11").Value = 0.6388539132363013
$ws.Range("P11").Value = 0.6388539132363011
$ws.Range("Q11").Value = 2666.284024000901
$ws.Range("R11").Value = 23996.55621600812
$ws.Range("S11").Value = 0.6377633460928189
$ws.Range("T11").Value = 0.6377633460928188
$ws.Range("G12").Value = 77.07094833333333
$ws.Range("H12").Value = 231.212845
$ws.Range("I12").Value = 0.9982929318880467
$ws.Range("J12").Value = 0.9982929318880466
$ws.Range("M12").Value = 6.728406666666667
$ws.Range("N12").Value = 20.18522
$ws.Range("O12").Value = 0.1242504825987572
$ws.Range("P12").Value = 0.1242504825987572
$ws.Range("Q12").Value = 518.5646825723222
$ws.Range("R12").Value = 4667.0821431509
$ws.Range("S12").Value = 0.1240383785620181
$ws.Range("T12").Value = 0.1240383785620181
$ws.Range("G13").Value = 77.07094833333333
$ws.Range("H13").Value = 231.212845
$ws.Range("I13").Value = 0.9982929318880467
$ws.Range("J13").Value = 0.9982929318880466
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 0.6666666666666666
$ws.Range("M13").Value = 0.1544473333333333
$ws.Range("N13").Value = 0.463342
$ws.Range("O13").Value = 0.002852109965027549
$ws.Range("P13").Value = 0.002852109965027548
$ws.Range("Q13").Value = 11.90340244755444
$ws.Range("R13").Value = 107.13062202799
$ws.Range("S13").Value = 0.002847241219054466
$ws.Range("T13").Value = 0.002847241219054465
